$wb = $excel.ActiveWorkbook
$zone42 = $wb.Worksheets.Item("zone42")
$zone42.Copy($zone42)
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "Sheet2"
$newSheet.Cells.ClearContents()
$newSheet.Range("A5:A1000").EntireRow.Delete()
$newSheet.Range("H1:T1").EntireColumn.Delete()
$newSheet.Cells.ClearFormats()
$newSheet.Columns("A:C").UseStandardWidth = $true
Write-Host "UsedRange:" $newSheet.UsedRange.Address(0,0)
